$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Marke:text; Farbe:text; Baujahr:number"
$ws.Range("D4").Value = "Tierart:text; Anzahl:number"
$ws.Range("D6").Value = "Land:text;Häufigkeit:select(Selten,Oft,Regelmäßig)"
$ws.Range("D7").Value = "Anteil:number;Technik:select(Laptop,PC,Tablet)"
$ws.Range("D9").Value = "Dringend(1 Woche):checkbox;Normal(2-4 Wochen):checkbox;Nicht-Dringend(5 Wochen und später:checkbox"

$ws.Range("D8").Select()
